$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header labels (shared strings) to the new dataset column names.
$ws.Range("B1").Value = "exportacion_bienes_servicios"
$ws.Range("C1").Value = "formacion_bruta_capital_fijo"
$ws.Range("D1").Value = "gasto_consumo_final_gobierno_central"
$ws.Range("E1").Value = "gasto_consumo_final_hogar_isflsh"
$ws.Range("F1").Value = "importaciones_bienes_servicios"
$ws.Range("G1").Value = "variacion_existencias"

# Swap the values of columns B and E for every data row (rows 2-62),
# since the two underlying series were mislabeled/swapped.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 62 }

for ($r = 2; $r -le $lastRow; $r++) {
    $bCell = $ws.Cells.Item($r, 2)
    $eCell = $ws.Cells.Item($r, 5)
    $bVal = $bCell.Value2
    $eVal = $eCell.Value2
    $bCell.Value2 = $eVal
    $eCell.Value2 = $bVal
}
